# Applies the cryptos.xlsx price/volume/coin-order update described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (coin names, links, percentage-change strings,
# and price strings that Excel would not mis-parse as numbers).
$simpleUpdates = @{
    'D2' = '68.407.61'
    'E2' = '  +2.06%  '
    'D3' = '2.642.53'
    'E3' = '  +1.50%  '
    'E4' = '  +0.02%  '
    'E5' = '  +1.64%  '
    'E6' = '  +3.12%  '
    'E7' = '  +0.00%  '
    'E8' = '  -0.26%  '
    'D9' = '2.641.91'
    'E9' = '  +1.55%  '
    'E10' = '  +5.46%  '
    'E11' = '  -0.44%  '
    'E12' = '  +1.44%  '
    'E13' = '  +1.48%  '
    'E14' = '  +2.43%  '
    'E15' = '  +2.59%  '
    'D16' = '3.123.24'
    'E16' = '  +1.59%  '
    'D17' = '68.277.94'
    'E17' = '  +2.11%  '
    'D18' = '2.664.77'
    'E18' = '  +2.55%  '
    'E19' = '  +3.26%  '
    'E20' = '  +0.33%  '
    'E21' = '  +0.63%  '
    'E22' = '  -0.66%  '
    'E23' = '  +2.50%  '
    'E24' = '  +2.63%  '
    'E25' = '  +0.45%  '
    'E26' = '  -0.01%  '
    'E27' = '  +1.49%  '
    'B28' = 'WrappedeETH'
    'C28' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D28' = '2.774.60'
    'E28' = '  +1.42%  '
    'B29' = 'PEPE'
    'C29' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'E29' = '  +5.81%  '
    'E30' = '  -0.43%  '
    'E31' = '  -1.13%  '
    'E32' = '  +4.88%  '
    'E33' = '  +4.48%  '
    'E34' = '  +2.42%  '
    'E35' = '  +3.32%  '
    'E36' = '  +0.03%  '
    'E37' = '  +3.67%  '
    'E38' = '  +2.53%  '
    'E39' = '  +4.31%  '
    'E40' = '  +1.22%  '
    'E41' = '  +3.23%  '
    'E42' = '  +0.47%  '
    'E43' = '  +2.42%  '
    'E44' = '  +3.79%  '
    'D45' = '0.0₆0321'
    'E45' = '  +12.47%  '
    'E46' = '  +0.05%  '
    'E47' = '  -0.12%  '
    'E48' = '  +2.63%  '
    'E49' = '  +0.65%  '
    'B50' = 'InjectiveProtocol'
    'C50' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'E50' = '  +2.33%  '
    'B51' = 'Optimism'
    'C51' = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
    'E51' = '  +1.38%  '
}
foreach ($cellRef in $simpleUpdates.Keys) {
    $ws.Range($cellRef).Value = $simpleUpdates[$cellRef]
}

# Price strings in column D that look like plain decimal numbers
# (e.g. '600.02', '1.00') must be forced to stay text, otherwise Excel's
# automatic type detection would silently convert them into numeric values
# and drop formatting such as trailing zeros. We write each value into a
# scratch cell pre-formatted as Text, copy it, and paste values+formats
# (PasteSpecial xlPasteAll = -4104) into the destination cell so the text
# type carries over without permanently attaching a new style to the
# scratch range. The helper column is removed afterwards.
$textUpdates = @{
    'D5' = '600.02'
    'D6' = '154.64'
    'D12' = '5.22'
    'D13' = '0.349'
    'D20' = '365.91'
    'D22' = '4.26'
    'D24' = '2.11'
    'D25' = '73.41'
    'D26' = '1.00'
    'D27' = '10.03'
    'D29' = '0.0000105'
    'D31' = '575.48'
    'D32' = '1.43'
    'D33' = '7.99'
    'D38' = '160.00'
    'D40' = '19.21'
    'D41' = '5.39'
    'D44' = '17.73'
    'D47' = '40.59'
    'D48' = '156.92'
    'D49' = '3.73'
    'D50' = '21.94'
    'D51' = '1.70'
}
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
foreach ($cellRef in $textUpdates.Keys) {
    $helper.Value = $textUpdates[$cellRef]
    $helper.Copy()
    $ws.Range($cellRef).PasteSpecial(-4104)
}
$ws.Columns("ZZ").Delete()
